$d = $word.ActiveDocument

function ReplaceText($find, $replace) {
    $r = $d.Content
    $ok = $r.Find.Execute($find, $false, $false, $false, $false, $false, $true, 1, $false, $replace, 2)
    if (-not $ok) {
        Write-Host "NOT FOUND: $find"
    }
    return $ok
}

function DeleteWholeParagraph($exactText) {
    $target = $null
    foreach ($p in $d.Paragraphs) {
        if ($p.Range.Text -eq ($exactText + "`r")) {
            $target = $p
            break
        }
    }
    if ($target -eq $null) {
        Write-Host "PARAGRAPH NOT FOUND: $exactText"
        return $false
    }
    $target.Range.Delete()
    return $true
}

# 1. "Wireframes site web." -> "Wireframes " (trailing space kept)
ReplaceText "Wireframes site web." "Wireframes "

# 2. Delete the whole "Wireframe logiciel." paragraph
DeleteWholeParagraph "Wireframe logiciel."

# 3. "Chemin utilisateur site web et logiciel." -> "Chemin utilisateur site web."
ReplaceText "Chemin utilisateur site web et logiciel." "Chemin utilisateur site web."

# 4. "Arborescence logiciel et web." -> "Arborescence."
ReplaceText "Arborescence logiciel et web." "Arborescence."

# 5. "Gestion des utilisateurs web, logiciel et réseau." -> "Gestion des utilisateurs."
ReplaceText "Gestion des utilisateurs web, logiciel et réseau." "Gestion des utilisateurs."

# 6. "Authentification utilisateur site web et logiciel" -> "UML" (leave trailing "." run untouched)
ReplaceText "Authentification utilisateur site web et logiciel" "UML"

# 7. Merge "Authentification Administrateur site web et logiciel." paragraph
#    with the following "Mise en forme et structuration des pages site web et logiciel." paragraph,
#    collapsing everything down to a single paragraph reading "Remplir doc."
#    Step 7a: collapse this paragraph's own text down to "Remplir doc"
$ok = ReplaceText "Authentification Administrateur site web et logiciel." "Remplir doc"
if ($ok) {
    # Step 7b: delete the paragraph mark right after "Remplir doc" so the following
    # paragraph ("Mise en forme et structuration des pages site web et logiciel.") merges into it
    $r2 = $d.Content
    $ok2 = $r2.Find.Execute("Remplir doc", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $ok2) {
        Write-Host "NOT FOUND (re-locate): Remplir doc"
    } else {
        $markRange = $d.Range($r2.End, $r2.End + 1)
        $markRange.Delete()

        # Step 7c: strip the now-merged-in leading text from the former second paragraph,
        # leaving only its trailing "." run attached to "Remplir doc"
        ReplaceText "Remplir docMise en forme et structuration des pages site web et logiciel" "Remplir doc"
    }
}

Write-Host "done"
